# Quarterly indexing esoteric bug-fix operation
#
# A new leading quarter value is inserted into column B of every data row
# (rows 2-16). This pushes each already-present quarter value one column to
# the right (B->C, C->D, ... J->K). Whatever previously sat in column K falls
# off the edge of the tracked series and is discarded.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: shift existing values right, then insert the new value at B2
# K2 is pushed past column K and is dropped
$ws.Range("K2").Value = $ws.Range("J2").Value2
$ws.Range("J2").Value = $ws.Range("I2").Value2
$ws.Range("I2").Value = $ws.Range("H2").Value2
$ws.Range("H2").Value = $ws.Range("G2").Value2
$ws.Range("G2").Value = $ws.Range("F2").Value2
$ws.Range("F2").Value = $ws.Range("E2").Value2
$ws.Range("E2").Value = $ws.Range("D2").Value2
$ws.Range("D2").Value = $ws.Range("C2").Value2
$ws.Range("C2").Value = $ws.Range("B2").Value2
$ws.Range("B2").Value = -0.2177157015159319

# Row 3: shift existing values right, then insert the new value at B3
# K3 is pushed past column K and is dropped
$ws.Range("K3").Value = $ws.Range("J3").Value2
$ws.Range("J3").Value = $ws.Range("I3").Value2
$ws.Range("I3").Value = $ws.Range("H3").Value2
$ws.Range("H3").Value = $ws.Range("G3").Value2
$ws.Range("G3").Value = $ws.Range("F3").Value2
$ws.Range("F3").Value = $ws.Range("E3").Value2
$ws.Range("E3").Value = $ws.Range("D3").Value2
$ws.Range("D3").Value = $ws.Range("C3").Value2
$ws.Range("C3").Value = $ws.Range("B3").Value2
$ws.Range("B3").Value = -0.1395947820665385

# Row 4: shift existing values right, then insert the new value at B4
# K4 is pushed past column K and is dropped
$ws.Range("K4").Value = $ws.Range("J4").Value2
$ws.Range("J4").Value = $ws.Range("I4").Value2
$ws.Range("I4").Value = $ws.Range("H4").Value2
$ws.Range("H4").Value = $ws.Range("G4").Value2
$ws.Range("G4").Value = $ws.Range("F4").Value2
$ws.Range("F4").Value = $ws.Range("E4").Value2
$ws.Range("E4").Value = $ws.Range("D4").Value2
$ws.Range("D4").Value = $ws.Range("C4").Value2
$ws.Range("C4").Value = $ws.Range("B4").Value2
$ws.Range("B4").Value = -0.3119065001142551

# Row 5: shift existing values right, then insert the new value at B5
# K5 is pushed past column K and is dropped
$ws.Range("K5").Value = $ws.Range("J5").Value2
$ws.Range("J5").Value = $ws.Range("I5").Value2
$ws.Range("I5").Value = $ws.Range("H5").Value2
$ws.Range("H5").Value = $ws.Range("G5").Value2
$ws.Range("G5").Value = $ws.Range("F5").Value2
$ws.Range("F5").Value = $ws.Range("E5").Value2
$ws.Range("E5").Value = $ws.Range("D5").Value2
$ws.Range("D5").Value = $ws.Range("C5").Value2
$ws.Range("C5").Value = $ws.Range("B5").Value2
$ws.Range("B5").Value = 0.7021231295320197

# Row 6: shift existing values right, then insert the new value at B6
# K6 is pushed past column K and is dropped
$ws.Range("K6").Value = $ws.Range("J6").Value2
$ws.Range("J6").Value = $ws.Range("I6").Value2
$ws.Range("I6").Value = $ws.Range("H6").Value2
$ws.Range("H6").Value = $ws.Range("G6").Value2
$ws.Range("G6").Value = $ws.Range("F6").Value2
$ws.Range("F6").Value = $ws.Range("E6").Value2
$ws.Range("E6").Value = $ws.Range("D6").Value2
$ws.Range("D6").Value = $ws.Range("C6").Value2
$ws.Range("C6").Value = $ws.Range("B6").Value2
$ws.Range("B6").Value = 1.514070997382048

# Row 7: shift existing values right, then insert the new value at B7
$ws.Range("K7").Value = $ws.Range("J7").Value2
$ws.Range("J7").Value = $ws.Range("I7").Value2
$ws.Range("I7").Value = $ws.Range("H7").Value2
$ws.Range("H7").Value = $ws.Range("G7").Value2
$ws.Range("G7").Value = $ws.Range("F7").Value2
$ws.Range("F7").Value = $ws.Range("E7").Value2
$ws.Range("E7").Value = $ws.Range("D7").Value2
$ws.Range("D7").Value = $ws.Range("C7").Value2
$ws.Range("C7").Value = $ws.Range("B7").Value2
$ws.Range("B7").Value = 0.2163102553365951

# Row 8: shift existing values right, then insert the new value at B8
$ws.Range("J8").Value = $ws.Range("I8").Value2
$ws.Range("I8").Value = $ws.Range("H8").Value2
$ws.Range("H8").Value = $ws.Range("G8").Value2
$ws.Range("G8").Value = $ws.Range("F8").Value2
$ws.Range("F8").Value = $ws.Range("E8").Value2
$ws.Range("E8").Value = $ws.Range("D8").Value2
$ws.Range("D8").Value = $ws.Range("C8").Value2
$ws.Range("C8").Value = $ws.Range("B8").Value2
$ws.Range("B8").Value = 0.3684555432821496

# Row 9: shift existing values right, then insert the new value at B9
$ws.Range("I9").Value = $ws.Range("H9").Value2
$ws.Range("H9").Value = $ws.Range("G9").Value2
$ws.Range("G9").Value = $ws.Range("F9").Value2
$ws.Range("F9").Value = $ws.Range("E9").Value2
$ws.Range("E9").Value = $ws.Range("D9").Value2
$ws.Range("D9").Value = $ws.Range("C9").Value2
$ws.Range("C9").Value = $ws.Range("B9").Value2
$ws.Range("B9").Value = 0.661541622456546

# Row 10: shift existing values right, then insert the new value at B10
$ws.Range("H10").Value = $ws.Range("G10").Value2
$ws.Range("G10").Value = $ws.Range("F10").Value2
$ws.Range("F10").Value = $ws.Range("E10").Value2
$ws.Range("E10").Value = $ws.Range("D10").Value2
$ws.Range("D10").Value = $ws.Range("C10").Value2
$ws.Range("C10").Value = $ws.Range("B10").Value2
$ws.Range("B10").Value = -0.07992401592518952

# Row 11: shift existing values right, then insert the new value at B11
$ws.Range("G11").Value = $ws.Range("F11").Value2
$ws.Range("F11").Value = $ws.Range("E11").Value2
$ws.Range("E11").Value = $ws.Range("D11").Value2
$ws.Range("D11").Value = $ws.Range("C11").Value2
$ws.Range("C11").Value = $ws.Range("B11").Value2
$ws.Range("B11").Value = 0.1551026493581833

# Row 12: shift existing values right, then insert the new value at B12
$ws.Range("F12").Value = $ws.Range("E12").Value2
$ws.Range("E12").Value = $ws.Range("D12").Value2
$ws.Range("D12").Value = $ws.Range("C12").Value2
$ws.Range("C12").Value = $ws.Range("B12").Value2
$ws.Range("B12").Value = -0.08373363042288225

# Row 13: shift existing values right, then insert the new value at B13
$ws.Range("E13").Value = $ws.Range("D13").Value2
$ws.Range("D13").Value = $ws.Range("C13").Value2
$ws.Range("C13").Value = $ws.Range("B13").Value2
$ws.Range("B13").Value = 0.1925427069667326

# Row 14: shift existing values right, then insert the new value at B14
$ws.Range("D14").Value = $ws.Range("C14").Value2
$ws.Range("C14").Value = $ws.Range("B14").Value2
$ws.Range("B14").Value = -0.4379379024501944

# Row 15: shift existing values right, then insert the new value at B15
$ws.Range("C15").Value = $ws.Range("B15").Value2
$ws.Range("B15").Value = 0.2324016585002178

# Row 16: shift existing values right, then insert the new value at B16
$ws.Range("B16").Value = -0.09587373626955231
